$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this data block (before row 969).
# Excel will shift all existing rows (969-1052) down by 2 (to 971-1054)
# and extend the used range / dimension automatically.
$ws.Rows("969:970").Insert()

# Populate the newly inserted row 969 (Primera) with the new weekly data point.
$ws.Cells.Item(969, 1).Value  = 3
$ws.Cells.Item(969, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(969, 3).Value  = "Coquimbo"
$ws.Cells.Item(969, 4).Value  = 45106
$ws.Cells.Item(969, 5).Value  = 5
$ws.Cells.Item(969, 6).Value  = "Fruta"
$ws.Cells.Item(969, 7).Value  = 100108
$ws.Cells.Item(969, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(969, 9).Value  = 100108005
$ws.Cells.Item(969, 10).Value = "Piña"
$ws.Cells.Item(969, 11).Value = "Caramelo"
$ws.Cells.Item(969, 12).Value = "Primera"
$ws.Cells.Item(969, 13).Value = 108
$ws.Cells.Item(969, 14).Value = 28000
$ws.Cells.Item(969, 15).Value = 28000
$ws.Cells.Item(969, 16).Value = 28000
$ws.Cells.Item(969, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(969, 18).Value = "Ecuador"
$ws.Cells.Item(969, 19).Value = 2333
$ws.Cells.Item(969, 20).Value = 12

# Populate the newly inserted row 970 (Segunda) with the new weekly data point.
$ws.Cells.Item(970, 1).Value  = 3
$ws.Cells.Item(970, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(970, 3).Value  = "Coquimbo"
$ws.Cells.Item(970, 4).Value  = 45106
$ws.Cells.Item(970, 5).Value  = 5
$ws.Cells.Item(970, 6).Value  = "Fruta"
$ws.Cells.Item(970, 7).Value  = 100108
$ws.Cells.Item(970, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(970, 9).Value  = 100108005
$ws.Cells.Item(970, 10).Value = "Piña"
$ws.Cells.Item(970, 11).Value = "Caramelo"
$ws.Cells.Item(970, 12).Value = "Segunda"
$ws.Cells.Item(970, 13).Value = 54
$ws.Cells.Item(970, 14).Value = 28000
$ws.Cells.Item(970, 15).Value = 28000
$ws.Cells.Item(970, 16).Value = 28000
$ws.Cells.Item(970, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(970, 18).Value = "Ecuador"
$ws.Cells.Item(970, 19).Value = 2000
$ws.Cells.Item(970, 20).Value = 14

# Make sure the date cells keep the same date-time number format as the
# rest of column D (style index 2 in the original workbook).
$ws.Range("D969:D970").NumberFormat = $ws.Range("D971").NumberFormat
